# parser: unit tests work
#
# Adds a fifth "all/f1_good/f1_bad/f2_good/f2_bad/f3_good/f3_bad/chi2_per_dof_th"
# summary block (rows 37-42, mirroring the header+labels pattern already used
# for the blocks at rows 1, 13, 21 and 29) and drops the now-unused
# "golay filter" tag that used to sit in K29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29's header used to carry an extra "golay filter" tag in K29 - remove it.
$ws.Range("K29").ClearContents() | Out-Null

# New header row 37 - same layout as rows 1 / 13 / 21 / 29.
$ws.Range("B37").Value = "all"
$ws.Range("C37").Value = "f1_good"
$ws.Range("D37").Value = "f1_bad"
$ws.Range("E37").Value = "f2_good"
$ws.Range("F37").Value = "f2_bad"
$ws.Range("G37").Value = "f3_good"
$ws.Range("H37").Value = "f3_bad"
$ws.Range("J37").Value = "chi2_per_dof_th"

# New data rows 38-42.
$ws.Range("A38").Value = "всего"
$ws.Range("J38").Value = 5

$ws.Range("A39").Value = "шумы"
$ws.Range("A40").Value = "одиночные"
$ws.Range("A41").Value = "двойные"
$ws.Range("A42").Value = "тройные"

# Match the author's final selection/cursor position.
$ws.Range("B38").Select() | Out-Null
